$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Fix the typo in the "Image de groupe 1" activity label (row 27, column B)
# before the row shift below so the text content ends up correct.
$ws.Range("B27").Value = "Image de groupe 1"

# Fill in the missing hours value for the "16/04/18" activity row.
$ws.Range("C17").Value = 2

# Remove the row for "21/5/18 - finalisation du rapport" entirely; Excel
# shifts all the rows below it up by one and keeps the SUM formula range
# self-adjusting.
$ws.Rows.Item(30).Delete()

# Update the view scroll position / selection to match the saved state.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("C26").Select()
